$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header labels (row 1), columns B..W
$headers = @{
    "B1" = "Investment"
    "C1" = "Saving"
    "D1" = "PROI"
    "E1" = "PPBT"
    "F1" = "Water Saving"
    "G1" = "Emission Saving"
    "H1" = "Land Saving"
    "I1" = "Import Saving"
    "J1" = "Capital Saving"
    "K1" = "Workforce Saving"
    "L1" = "Water Investment"
    "M1" = "Emission Investment"
    "N1" = "Land Investment"
    "O1" = "Import Investment"
    "P1" = "Workforce Investment"
    "Q1" = "Capital Investment"
    "R1" = "Water Total Impact"
    "S1" = "Emission Total Impact"
    "T1" = "Land Total Impact"
    "U1" = "Import Total Impact"
    "V1" = "Workforce Total Impact"
    "W1" = "Capital Total Impact"
}

foreach ($addr in $headers.Keys) {
    $ws.Range($addr).Value = $headers[$addr]
}

# New data values (row 2), columns B..W
$values = @{
    "B2" = 1.000000158324838
    "C2" = 0.9354393007233739
    "D2" = 0.9354391526201219
    "E2" = 1.06901661877103
    "F2" = 0.04106397050782107
    "G2" = 0.05008577914486523
    "H2" = 0.003466384087914776
    "I2" = 0.08160151727497578
    "J2" = 0.07793689612299204
    "K2" = 0.06345285149291158
    "L2" = 0.0005181713204365224
    "M2" = 0.0009185673316096654
    "N2" = 0.000001404879185429309
    "O2" = 0.06491394690237939
    "P2" = 0.09720575390383601
    "Q2" = 0.00409148377366364
    "R2" = -0.8207612388359848
    "S2" = -1.000797015565695
    "T2" = -0.06932627687911008
    "U2" = -1.567116398597136
    "V2" = -1.171851275954396
    "W2" = -1.554646438686177
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
